$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Testing password from the file" -> append completion sentence.
#    The unfinished sentence gets its ending typed in by Jordan.
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Testing password from the file", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$r1.Collapse(0)
$r1.InsertAfter(" will be completed by Jordan")

# ------------------------------------------------------------------
# 2. Finish the "Riley wants to create GUI..." sentence:
#    ", asked jeff to populate the list." -> " and asked Jeffrey to populate the list."
#    Locate the whole sentence first to get its exact bounds, then
#    build a freshly-scoped Range over just that sentence so the
#    replace (and later the bookmark placement) cannot touch the
#    other "jeff"/"Jeffrey" occurrences elsewhere in the document.
# ------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("Riley wants to create GUI for view enrollments file, asked jeff to populate the list.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentenceStart = $find2.Start
$sentenceEnd = $find2.End

$scoped2 = $d.Range($sentenceStart, $sentenceEnd)
$scoped2.Find.Execute("file, asked jeff to populate the list.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "file and asked Jeffrey to populate the list.", 2)

# ------------------------------------------------------------------
# 3. Move the _GoBack bookmark so it sits right after "Jeffrey" in
#    the sentence above instead of in the trailing empty paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$find3 = $d.Content
$find3.Find.Execute("Riley wants to create GUI for view enrollments file and asked Jeffrey to populate the list.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newSentenceStart = $find3.Start
$newSentenceEnd = $find3.End

$scoped3 = $d.Range($newSentenceStart, $newSentenceEnd)
$scoped3.Find.Execute("Jeffrey", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scoped3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $scoped3)
